# Atualizado por script em 12-11-2023 14:45
#
# This script applies the betexplorer-scrape update to the Poland
# division-2 2023-2024 sheet:
#  1) Four pairs of adjacent match rows had their match-detail columns
#     (F..V) swapped (the A..E "index/meta" columns stay put) -- this
#     reflects the two fixtures on the same match-day being re-ordered
#     by the scraper.
#  2) Four brand-new match rows (147-150) were appended at the bottom,
#     for matches played 11/11/2023, scraped 12/11/2023.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowDetails {
    param($rowA, $rowB)

    $rangeA = $ws.Range("F$($rowA):V$($rowA)")
    $rangeB = $ws.Range("F$($rowB):V$($rowB)")

    $valsA = $rangeA.Value2
    $valsB = $rangeB.Value2

    $rangeA.Value2 = $valsB
    $rangeB.Value2 = $valsA
}

# --- 1) Swap F:V between the described adjacent row pairs -------------
Swap-RowDetails 6 7
Swap-RowDetails 92 93
Swap-RowDetails 103 104
Swap-RowDetails 136 137

# --- 2) Append four new rows (147-150), copying formatting from the --
#        last existing data row (146) before writing the new values.
$ws.Range("A146:V146").Copy()
$ws.Range("A147:V150").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = 0

$newRows = @(
    @{ Row=147; A=146; B="poland"; C="division-2"; D="2023-2024"; E=45242.52083333334;
       F="GKS Jastrzebie"; G=3; H="LKS Lodz II"; I=2;
       J=2.15; K="11/11/2023 00:43"; L=2.15; M="12/11/2023 12:28";
       N=3.31; O="11/11/2023 00:43"; P=3.71; Q="12/11/2023 12:28";
       R=2.88; S="11/11/2023 00:43"; T=2.98; U="12/11/2023 12:28";
       V="https://www.betexplorer.com/football/poland/division-2/gks-jastrzebie-lks-lodz/I7L0WUQR/" },
    @{ Row=148; A=147; B="poland"; C="division-2"; D="2023-2024"; E=45242.54166666666;
       F="Lech Poznan II"; G=1; H="Polonia Bytom"; I=0;
       J=2.89; K="11/11/2023 01:13"; L=2.98; M="12/11/2023 12:51";
       N=3.36; O="11/11/2023 01:13"; P=3.55; Q="12/11/2023 12:51";
       R=2.17; S="11/11/2023 01:13"; T=2.21; U="12/11/2023 12:51";
       V="https://www.betexplorer.com/football/poland/division-2/lech-poznan-polonia-bytom/8IybC9Ze/" },
    @{ Row=149; A=148; B="poland"; C="division-2"; D="2023-2024"; E=45242.54166666666;
       F="Stezyca"; G=1; H="Stomil Olsztyn"; I=0;
       J=1.95; K="11/11/2023 01:13"; L=1.9; M="12/11/2023 12:39";
       N=3.28; O="11/11/2023 01:13"; P=3.33; Q="12/11/2023 12:39";
       R=3.35; S="11/11/2023 01:13"; T=4.08; U="12/11/2023 12:39";
       V="https://www.betexplorer.com/football/poland/division-2/stezyca-stomil-olsztyn/dhXiYjdF/" },
    @{ Row=150; A=149; B="poland"; C="division-2"; D="2023-2024"; E=45242.54166666666;
       F="Zaglebie II"; G=4; H="S. Wola"; I=0;
       J=2.3; K="11/11/2023 01:13"; L=2.72; M="12/11/2023 12:51";
       N=3.19; O="11/11/2023 01:13"; P=3.15; Q="12/11/2023 12:51";
       R=2.81; S="11/11/2023 01:13"; T=2.59; U="12/11/2023 12:51";
       V="https://www.betexplorer.com/football/poland/division-2/zaglebie-stal-stalowa-wola/zqWeXABL/" }
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

foreach ($rowData in $newRows) {
    $r = $rowData.Row
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value = $rowData[$col]
    }
}
